$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.727.70"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "2.095.58"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5164"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4387"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09257"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "2.102.18"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.775"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.235"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06660"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "29.764.20"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.318"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "2.347.40"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.520"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.140"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.00%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.641"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.177"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.323"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02577"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06731"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.323"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2228"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7012"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.322"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.627"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000356"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.86%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
